# Bond screener "as of" date roll-forward by exactly one day
# (2023-10-15 serial 45214 -> 2023-10-16 serial 45215).
#
# Columns on the data rows (row 1 is the header):
#   F = date of the previous coupon payment (date-formatted serial number)
#   G = number of days elapsed since that previous payment  (= asOf - F)
#   H = date of the next coupon payment      (date-formatted serial number)
#   I = number of days remaining until that next payment     (= H - asOf)
#
# Moving "today" forward by one day means:
#   - every "days since" counter (G) goes up by 1
#   - every "days until" counter (I) goes down by 1
#   - UNLESS the previous payment was exactly one full coupon period
#     (183 days) before the old "as of" date, in which case a new coupon
#     just started: F becomes the old "as of" date (45214) and G resets to 1
#     (H/I just keep counting down normally).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$asOfBefore = 45214
$rolloverPeriod = 183

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)   # F: previous payout date
    $gCell = $ws.Cells.Item($r, 7)   # G: days since previous payout
    $iCell = $ws.Cells.Item($r, 9)   # I: days until next payout

    $fVal = $fCell.Value2()
    $gVal = $gCell.Value2()
    $iVal = $iCell.Value2()

    if ($null -ne $fVal -and $null -ne $gVal) {
        if ($fVal -eq ($asOfBefore - $rolloverPeriod)) {
            # The coupon period just rolled over: a new payment happened
            # exactly on the old "as of" date.
            $fCell.Value = $asOfBefore
            $gCell.Value = 1
        } else {
            $gCell.Value = $gVal + 1
        }
    }

    if ($null -ne $iVal) {
        $iCell.Value = $iVal - 1
    }
}

Write-Output "Rolled bond screener dates forward by one day"
